$d = $word.ActiveDocument

$replacements = @(
    @{old="743×5="; new="275×6="},
    @{old="541×7="; new="543×3="},
    @{old="944×7="; new="724×9="},
    @{old="892×8="; new="847×2="},
    @{old="190×2="; new="439×4="},
    @{old="721×9="; new="173×9="},
    @{old="843×9="; new="417×5="},
    @{old="331×7="; new="318×6="},
    @{old="365×2="; new="806×4="},
    @{old="127×2="; new="984×6="},
    @{old="982×2="; new="320×6="},
    @{old="535×3="; new="886×6="},
    @{old="616×4="; new="388×4="},
    @{old="612×2="; new="378×4="},
    @{old="635×3="; new="144×2="},
    @{old="847×4="; new="595×8="},
    @{old="388×6="; new="889×6="},
    @{old="848×5="; new="845×6="},
    @{old="251×5="; new="551×5="},
    @{old="930×9="; new="494×7="},
    @{old="701×6="; new="110×6="},
    @{old="497×7="; new="437×5="},
    @{old="262×7="; new="598×3="},
    @{old="500×7="; new="198×9="},
    @{old="834×3="; new="144×3="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
